$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 6541.5884
$ws.Range("I41").Value = 782.1667
$ws.Range("J41").Value = 9683.091
$ws.Range("K41").Value = 782.1667
$ws.Range("L41").Value = 9683.091
$ws.Range("M41").Value = -342.1667
$ws.Range("N41").Value = -10563.091
$ws.Range("H51").Value = 3429.7307
$ws.Range("I51").Value = 3918.3635
$ws.Range("J51").Value = 3071.4
$ws.Range("K51").Value = 3918.3635
$ws.Range("L51").Value = 3071.4
$ws.Range("M51").Value = -3434.3635
$ws.Range("N51").Value = -4039.4
$ws.Range("H53").Value = 198.84
$ws.Range("I53").Value = 275.33334
$ws.Range("J53").Value = 155.8125
$ws.Range("K53").Value = 275.33334
$ws.Range("L53").Value = 155.8125
$ws.Range("M53").Value = 361.66666
$ws.Range("N53").Value = -1429.8125
$ws.Range("H62").Value = 11965.654
$ws.Range("I62").Value = 25611.223
$ws.Range("J62").Value = 4741.5293
$ws.Range("K62").Value = 25611.223
$ws.Range("L62").Value = 4741.5293
$ws.Range("M62").Value = -24987.223
$ws.Range("N62").Value = -5989.5293
$ws.Range("H65").Value = 11965.654
$ws.Range("I65").Value = 25611.223
$ws.Range("J65").Value = 4741.5293
$ws.Range("K65").Value = 128056.115
$ws.Range("L65").Value = 23707.6465
$ws.Range("M65").Value = -124936.115
$ws.Range("N65").Value = -29947.6465
$ws.Range("H116").Value = 4132.4443
$ws.Range("I116").Value = 3886.5715
$ws.Range("J116").Value = 4993
$ws.Range("K116").Value = 3886.5715
$ws.Range("L116").Value = 4993
$ws.Range("M116").Value = -444.5715
$ws.Range("N116").Value = -11877
$ws.Range("H131").Value = 2373.1943
$ws.Range("I131").Value = 1146.5625
$ws.Range("J131").Value = 3354.5
$ws.Range("K131").Value = 3439.6875
$ws.Range("L131").Value = 10063.5
$ws.Range("M131").Value = 1600.3125
$ws.Range("N131").Value = -20143.5
$ws.Range("H132").Value = 3759.5908
$ws.Range("I132").Value = 3870
$ws.Range("J132").Value = 3600.111
$ws.Range("K132").Value = 11610
$ws.Range("L132").Value = 10800.333
$ws.Range("M132").Value = -9080
$ws.Range("N132").Value = -15860.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1036.6666
$ws.Range("I2").Value = 793.75
$ws.Range("J2").Value = 2980
$ws.Range("K2").Value = 793.75
$ws.Range("L2").Value = 2980
$ws.Range("M2").Value = -680.75
$ws.Range("N2").Value = -3206
$ws.Range("H45").Value = 1006.75
$ws.Range("I45").Value = 929.6667
$ws.Range("J45").Value = 1238
$ws.Range("K45").Value = 929.6667
$ws.Range("L45").Value = 1238
$ws.Range("M45").Value = -552.6667
$ws.Range("N45").Value = -1992
$ws.Range("H88").Value = 2558.862
$ws.Range("I88").Value = 2556
$ws.Range("K88").Value = 2556
$ws.Range("M88").Value = -2150
$ws.Range("H91").Value = 2558.862
$ws.Range("I91").Value = 2556
$ws.Range("K91").Value = 2556
$ws.Range("M91").Value = -1152
$ws.Range("H110").Value = 770.5714
$ws.Range("I110").Value = 761.5
$ws.Range("J110").Value = 825
$ws.Range("K110").Value = 761.5
$ws.Range("L110").Value = 825
$ws.Range("M110").Value = 1283.5
$ws.Range("N110").Value = -4915
$ws.Range("H116").Value = 1036.6666
$ws.Range("I116").Value = 793.75
$ws.Range("J116").Value = 2980
$ws.Range("K116").Value = 793.75
$ws.Range("L116").Value = 2980
$ws.Range("M116").Value = 1500.25
$ws.Range("N116").Value = -7568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1036.6666
$ws.Range("I3").Value = 793.75
$ws.Range("J3").Value = 2980
$ws.Range("K3").Value = 793.75
$ws.Range("L3").Value = 2980
$ws.Range("M3").Value = -679.75
$ws.Range("N3").Value = -3208
$ws.Range("H99").Value = 6619
$ws.Range("I99").Value = 7132.222
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 7132.222
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -5634.222
$ws.Range("N99").Value = -4996
$ws.Range("H134").Value = 5097.0625
$ws.Range("I134").Value = 5267.3423
$ws.Range("K134").Value = 15802.0269
$ws.Range("M134").Value = -13267.0269

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2873.1155
$ws.Range("I31").Value = 2088.4
$ws.Range("J31").Value = 3943.182
$ws.Range("K31").Value = 2088.4
$ws.Range("L31").Value = 3943.182
$ws.Range("M31").Value = -1793.4
$ws.Range("N31").Value = -4533.182
$ws.Range("H34").Value = 2873.1155
$ws.Range("I34").Value = 2088.4
$ws.Range("J34").Value = 3943.182
$ws.Range("K34").Value = 2088.4
$ws.Range("L34").Value = 3943.182
$ws.Range("M34").Value = -1886.4
$ws.Range("N34").Value = -4347.182
$ws.Range("H58").Value = 3839.275
$ws.Range("I58").Value = 4213.4644
$ws.Range("J58").Value = 2966.1667
$ws.Range("K58").Value = 4213.4644
$ws.Range("L58").Value = 2966.1667
$ws.Range("M58").Value = -4010.4644
$ws.Range("N58").Value = -3372.1667
$ws.Range("H99").Value = 86453.914
$ws.Range("I99").Value = 202082.4
$ws.Range("J99").Value = 3862.1428
$ws.Range("K99").Value = 202082.4
$ws.Range("L99").Value = 3862.1428
$ws.Range("M99").Value = -200584.4
$ws.Range("N99").Value = -6858.1428
$ws.Range("H126").Value = 86453.914
$ws.Range("I126").Value = 202082.4
$ws.Range("J126").Value = 3862.1428
$ws.Range("K126").Value = 606247.2
$ws.Range("L126").Value = 11586.4284
$ws.Range("M126").Value = -603777.2
$ws.Range("N126").Value = -16526.4284
$ws.Range("H134").Value = 1582.3611
$ws.Range("I134").Value = 995.03845
$ws.Range("K134").Value = 2985.11535
$ws.Range("M134").Value = -450.11535
$ws.Range("H136").Value = 3839.275
$ws.Range("I136").Value = 4213.4644
$ws.Range("J136").Value = 2966.1667
$ws.Range("K136").Value = 12640.3932
$ws.Range("L136").Value = 8898.500100000001
$ws.Range("M136").Value = -10090.3932
$ws.Range("N136").Value = -13998.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4112.346
$ws.Range("I126").Value = 4159.857
$ws.Range("J126").Value = 4056.9167
$ws.Range("K126").Value = 12479.571
$ws.Range("L126").Value = 12170.7501
$ws.Range("M126").Value = -10009.571
$ws.Range("N126").Value = -17110.7501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2772
$ws.Range("I122").Value = 2727.8462
$ws.Range("J122").Value = 2963.3333
$ws.Range("K122").Value = 8183.5386
$ws.Range("L122").Value = 8889.999899999999
$ws.Range("M122").Value = -5733.5386
$ws.Range("N122").Value = -13789.9999
$ws.Range("H136").Value = 3369.475
$ws.Range("J136").Value = 3647.7778
$ws.Range("L136").Value = 10943.3334
$ws.Range("N136").Value = -16043.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 404
$ws.Range("I113").Value = 288.66666
$ws.Range("J113").Value = 750
$ws.Range("K113").Value = 865.9999799999999
$ws.Range("L113").Value = 2250
$ws.Range("M113").Value = 1304.00002
$ws.Range("N113").Value = -6590
$ws.Range("H132").Value = 2137.303
$ws.Range("I132").Value = 990.2222
$ws.Range("K132").Value = 2970.6666
$ws.Range("M132").Value = -440.6666
